$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Modelo" header column (F1), matching the style of the existing
# header cells (e.g. E1: bold, bordered, centered).
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Add the model name value for the data row (F2), no special styling,
# matching the plain data cells like E2.
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"

$excel.CutCopyMode = $false
